# Update the cached "Last updated" date shown by the auto-updating
# Date placeholder (datetime / datetimeFigureOut field) from 8/31/22 to
# 9/2/22 everywhere it is cached: the Slide Master, all 14 slide
# layouts, and the Notes Master.

$p = $ppt.ActivePresentation
$oldDate = "8/31/22"
$newDate = "9/2/22"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.HasTextFrame) {
                if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                    $sh.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# Slide Master
Update-DateShape $p.SlideMaster.Shapes

# Every slide layout (CustomLayout) under the slide master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $cl = $layouts.Item($li)
    Update-DateShape $cl.Shapes
}

# Notes Master
Update-DateShape $p.NotesMaster.Shapes
